$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 302.26666
$ws.Range("I28").Value = 367
$ws.Range("K28").Value = 367
$ws.Range("M28").Value = 118
$ws.Range("H62").Value = 2183.6875
$ws.Range("I62").Value = 1983.5834
$ws.Range("K62").Value = 1983.5834
$ws.Range("M62").Value = -1359.5834
$ws.Range("H65").Value = 2183.6875
$ws.Range("I65").Value = 1983.5834
$ws.Range("K65").Value = 9917.916999999999
$ws.Range("M65").Value = -6797.916999999999
$ws.Range("H100").Value = 2950.625
$ws.Range("I100").Value = 2826.25
$ws.Range("J100").Value = 3075
$ws.Range("K100").Value = 2826.25
$ws.Range("L100").Value = 3075
$ws.Range("M100").Value = -2285.25
$ws.Range("N100").Value = -4157
$ws.Range("H106").Value = 1367.909
$ws.Range("I106").Value = 1100.4445
$ws.Range("K106").Value = 1100.4445
$ws.Range("M106").Value = -469.4445000000001
$ws.Range("H135").Value = 27786886
$ws.Range("I135").Value = 1149.8462
$ws.Range("K135").Value = 10348.6158
$ws.Range("M135").Value = -7813.6158
$ws.Range("H137").Value = 30331.777
$ws.Range("I137").Value = 2962.6316
$ws.Range("J137").Value = 60920.824
$ws.Range("K137").Value = 8887.8948
$ws.Range("L137").Value = 182762.472
$ws.Range("M137").Value = -6337.8948
$ws.Range("N137").Value = -187862.472
$ws.Range("H138").Value = 2490.75
$ws.Range("I138").Value = 741.6667
$ws.Range("J138").Value = 3285.7878
$ws.Range("K138").Value = 2225.0001
$ws.Range("L138").Value = 9857.3634
$ws.Range("M138").Value = 2914.9999
$ws.Range("N138").Value = -20137.3634
$ws.Range("H141").Value = 1774.8529
$ws.Range("I141").Value = 1197.9032
$ws.Range("J141").Value = 7736.6665
$ws.Range("K141").Value = 3593.7096
$ws.Range("L141").Value = 23209.9995
$ws.Range("M141").Value = 1586.2904
$ws.Range("N141").Value = -33569.99950000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 961.11475
$ws.Range("I2").Value = 740.44446
$ws.Range("J2").Value = 1581.75
$ws.Range("K2").Value = 740.44446
$ws.Range("L2").Value = 1581.75
$ws.Range("M2").Value = -627.44446
$ws.Range("N2").Value = -1807.75
$ws.Range("H32").Value = 16753.646
$ws.Range("I32").Value = 18947.05
$ws.Range("K32").Value = 18947.05
$ws.Range("M32").Value = -18660.05
$ws.Range("H61").Value = 490050.22
$ws.Range("I61").Value = 722685.2
$ws.Range("J61").Value = 5394
$ws.Range("K61").Value = 722685.2
$ws.Range("L61").Value = 5394
$ws.Range("M61").Value = -722473.2
$ws.Range("N61").Value = -5818
$ws.Range("H97").Value = 2623.375
$ws.Range("I97").Value = 2334.5
$ws.Range("J97").Value = 3490
$ws.Range("K97").Value = 2334.5
$ws.Range("L97").Value = 3490
$ws.Range("M97").Value = -1838.5
$ws.Range("N97").Value = -4482
$ws.Range("H116").Value = 961.11475
$ws.Range("I116").Value = 740.44446
$ws.Range("J116").Value = 1581.75
$ws.Range("K116").Value = 740.44446
$ws.Range("L116").Value = 1581.75
$ws.Range("M116").Value = 1553.55554
$ws.Range("N116").Value = -6169.75
$ws.Range("H122").Value = 2383.606
$ws.Range("I122").Value = 2281.84
$ws.Range("K122").Value = 6845.52
$ws.Range("M122").Value = -4395.52
$ws.Range("H136").Value = 490050.22
$ws.Range("I136").Value = 722685.2
$ws.Range("J136").Value = 5394
$ws.Range("K136").Value = 2168055.6
$ws.Range("L136").Value = 16182
$ws.Range("M136").Value = -2165505.6
$ws.Range("N136").Value = -21282

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 961.11475
$ws.Range("I3").Value = 740.44446
$ws.Range("J3").Value = 1581.75
$ws.Range("K3").Value = 740.44446
$ws.Range("L3").Value = 1581.75
$ws.Range("M3").Value = -626.44446
$ws.Range("N3").Value = -1809.75
$ws.Range("H99").Value = 2018.1666
$ws.Range("I99").Value = 2018.1666
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2018.1666
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -520.1666
$ws.Range("N99").Value = $null
$ws.Range("H107").Value = 625
$ws.Range("I107").Value = 625
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 625
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1295
$ws.Range("N107").Value = $null
$ws.Range("H134").Value = 63272.06
$ws.Range("I134").Value = 63272.06
$ws.Range("K134").Value = 189816.18
$ws.Range("M134").Value = -187281.18

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9830.333000000001
$ws.Range("I31").Value = 22720.066
$ws.Range("J31").Value = 2669.3704
$ws.Range("K31").Value = 22720.066
$ws.Range("L31").Value = 2669.3704
$ws.Range("M31").Value = -22425.066
$ws.Range("N31").Value = -3259.3704
$ws.Range("H34").Value = 9830.333000000001
$ws.Range("I34").Value = 22720.066
$ws.Range("J34").Value = 2669.3704
$ws.Range("K34").Value = 22720.066
$ws.Range("L34").Value = 2669.3704
$ws.Range("M34").Value = -22518.066
$ws.Range("N34").Value = -3073.3704
$ws.Range("H122").Value = 2025
$ws.Range("I122").Value = 2418.5
$ws.Range("K122").Value = 7255.5
$ws.Range("M122").Value = -4805.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 804.8
$ws.Range("I5").Value = 650.0769
$ws.Range("J5").Value = 972.4167
$ws.Range("K5").Value = 1950.2307
$ws.Range("L5").Value = 2917.2501
$ws.Range("M5").Value = -1838.2307
$ws.Range("N5").Value = -3141.2501
$ws.Range("H68").Value = 3637.8206
$ws.Range("I68").Value = 839.8
$ws.Range("J68").Value = 4049.2942
$ws.Range("K68").Value = 2519.4
$ws.Range("L68").Value = 12147.8826
$ws.Range("M68").Value = -1708.4
$ws.Range("N68").Value = -13769.8826
$ws.Range("H71").Value = 3637.8206
$ws.Range("I71").Value = 839.8
$ws.Range("J71").Value = 4049.2942
$ws.Range("K71").Value = 7558.2
$ws.Range("L71").Value = 36443.6478
$ws.Range("M71").Value = -3502.2
$ws.Range("N71").Value = -44555.6478
$ws.Range("H107").Value = 5129.4585
$ws.Range("J107").Value = 1057.4762
$ws.Range("L107").Value = 3172.4286
$ws.Range("N107").Value = -7012.4286
$ws.Range("H113").Value = 14897.143
$ws.Range("I113").Value = 33800
$ws.Range("J113").Value = 720
$ws.Range("K113").Value = 101400
$ws.Range("L113").Value = 2160
$ws.Range("M113").Value = -99230
$ws.Range("N113").Value = -6500
$ws.Range("H123").Value = 3698.5715
$ws.Range("J123").Value = 4115
$ws.Range("L123").Value = 12345
$ws.Range("N123").Value = -17245
$ws.Range("H131").Value = 108364.31
$ws.Range("I131").Value = 794.2857
$ws.Range("J131").Value = 117120.01
$ws.Range("K131").Value = 2382.8571
$ws.Range("L131").Value = 351360.03
$ws.Range("M131").Value = 2657.1429
$ws.Range("N131").Value = -361440.03
$ws.Range("H135").Value = 804.8
$ws.Range("I135").Value = 650.0769
$ws.Range("J135").Value = 972.4167
$ws.Range("K135").Value = 5850.6921
$ws.Range("L135").Value = 8751.7503
$ws.Range("M135").Value = -3315.6921
$ws.Range("N135").Value = -13821.7503
$ws.Range("H140").Value = 1851.7646
$ws.Range("I140").Value = 1468
$ws.Range("K140").Value = 4404
$ws.Range("M140").Value = 776
$ws.Range("H141").Value = 1448
$ws.Range("I141").Value = 419.66666
$ws.Range("J141").Value = 4533
$ws.Range("K141").Value = 1258.99998
$ws.Range("L141").Value = 13599
$ws.Range("M141").Value = 3921.00002
$ws.Range("N141").Value = -23959

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4582.32
$ws.Range("I7").Value = 4717.7
$ws.Range("J7").Value = 4040.8
$ws.Range("K7").Value = 4717.7
$ws.Range("L7").Value = 4040.8
$ws.Range("M7").Value = -4605.7
$ws.Range("N7").Value = -4264.8
$ws.Range("H22").Value = 2416.6
$ws.Range("I22").Value = 2800.25
$ws.Range("K22").Value = 2800.25
$ws.Range("M22").Value = -2505.25
$ws.Range("H27").Value = 2416.6
$ws.Range("I27").Value = 2800.25
$ws.Range("K27").Value = 2800.25
$ws.Range("M27").Value = -2693.25
$ws.Range("H93").Value = 1848.05
$ws.Range("I93").Value = 1960.6471
$ws.Range("K93").Value = 1960.6471
$ws.Range("M93").Value = -712.6470999999999
$ws.Range("H122").Value = 3143.647
$ws.Range("I122").Value = 2413.7
$ws.Range("K122").Value = 7241.099999999999
$ws.Range("M122").Value = -4791.099999999999
$ws.Range("H124").Value = 35424.5
$ws.Range("J124").Value = 35424.5
$ws.Range("L124").Value = 35424.5
$ws.Range("N124").Value = -45244.5
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H126").Value = 4582.32
$ws.Range("I126").Value = 4717.7
$ws.Range("J126").Value = 4040.8
$ws.Range("K126").Value = 14153.1
$ws.Range("L126").Value = 12122.4
$ws.Range("M126").Value = -11683.1
$ws.Range("N126").Value = -17062.4
$ws.Range("H132").Value = 1508.3077
$ws.Range("I132").Value = 1269.826
$ws.Range("J132").Value = 3336.6667
$ws.Range("K132").Value = 3809.478
$ws.Range("L132").Value = 10010.0001
$ws.Range("M132").Value = -1279.478
$ws.Range("N132").Value = -15070.0001
$ws.Range("H136").Value = 2073.762
$ws.Range("I136").Value = 1599.6428
$ws.Range("J136").Value = 3022
$ws.Range("K136").Value = 4798.928400000001
$ws.Range("L136").Value = 9066
$ws.Range("M136").Value = -2248.928400000001
$ws.Range("N136").Value = -14166

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14135.714
$ws.Range("J54").Value = 14135.714
$ws.Range("L54").Value = 14135.714
$ws.Range("N54").Value = -15175.714
$ws.Range("H122").Value = 1577.7084
$ws.Range("I122").Value = 1421.05
$ws.Range("K122").Value = 4263.15
$ws.Range("M122").Value = -1813.15
$ws.Range("H132").Value = 3187.5
$ws.Range("I132").Value = 3142.8572
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 9428.571599999999
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -6898.571599999999
$ws.Range("N132").Value = -15560
